$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H7").Value = '2x2 matrix'
$ws.Range("H10").Value = 'conceptual diagram'
$ws.Range("H12").Value = 'conceptual diagram'
$ws.Range("H13").Value = 'conceptual diagram'
$ws.Range("H14").Value = 'conceptual diagram'
$ws.Range("H15").Value = 'process diagram'
$ws.Range("H16").Value = 'conceptual diagram'
$ws.Range("H17").Value = 'photo'
$ws.Range("H18").Value = 'conceptual diagram'
$ws.Range("H19").Value = 'conceptual diagram'
$ws.Range("H22").Value = 'conceptual diagram'
$ws.Range("H24").Value = 'cycle'
$ws.Range("H25").Value = 'cycle'
$ws.Range("H26").Value = 'conceptual diagram'
$ws.Range("H30").Value = 'mixed statistical plot (more than 1 statistical plot type)'
$ws.Range("H31").Value = 'conceptual diagram'
$ws.Range("H34").Value = 'conceptual diagram'
$ws.Range("H35").Value = 'cycle'
$ws.Range("H37").Value = 'conceptual diagram'
$ws.Range("H38").Value = 'conceptual diagram'
$ws.Range("H51").Value = 'process diagram'
$ws.Range("H52").Value = 'conceptual diagram'
$ws.Range("H53").Value = 'conceptual diagram'
$ws.Range("H56").Value = 'process diagram'
$ws.Range("H57").Value = 'photo'
$ws.Range("H58").Value = 'conceptual diagram'
$ws.Range("H60").Value = 'photo'
$ws.Range("H62").Value = 'conceptual diagram'
$ws.Range("H63").Value = 'conceptual diagram'
$ws.Range("H64").Value = 'conceptual diagram'
